$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally has 8 rows (1 header + 7 full vaccine-name rows).
# For each full-name row, a new row with the abbreviated vaccine name is
# inserted directly beneath it (Johnson & Johnson / Janssen gets two new
# rows). First build the 16-row skeleton by inserting blank rows (working
# bottom-up so earlier row numbers stay stable), then fill in the new
# abbreviated names.

# After row 8 (Covishield ... ) -> row 9 is already empty once other rows
# below have shifted; nothing to insert for the very last row.

# After row 7 (Sinopharm BBIBP-CorV Vero Cells) -> insert new row 8.
$ws.Rows.Item(8).Insert()

# After row 6 (Oxford/AstraZeneca AZD1222) -> insert new row 7.
$ws.Rows.Item(7).Insert()

# After row 5 (Sinovac CoronaVac) -> insert new row 6.
$ws.Rows.Item(6).Insert()

# After row 4 (Johnson & Johnson Janssen Ad26.COV2.S) -> insert two new rows.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(6).Insert()

# After row 3 (Pfizer/BioNTech BNT162b2) -> insert new row 4.
$ws.Rows.Item(4).Insert()

# After row 2 (Moderna mRNA-1273) -> insert new row 3.
$ws.Rows.Item(3).Insert()

# Now the 16-row skeleton is in place. Fill the new rows with the
# abbreviated vaccine names.
$ws.Range("A5").Value = "Pfizer"
$ws.Range("A3").Value = "Moderna"
$ws.Range("A7").Value = "Johnson & Johnson"
$ws.Range("A8").Value = "Janssen"
$ws.Range("A12").Value = "AstraZeneca"
$ws.Range("A16").Value = "Covishield"
$ws.Range("A14").Value = "Sinopharm"
$ws.Range("A10").Value = "Sinovac"

# Clear formatting on the new rows so they don't inherit the bold/category
# style of the row above.
$ws.Range("A3").ClearFormats()
$ws.Range("A5").ClearFormats()
$ws.Range("A7").ClearFormats()
$ws.Range("A8").ClearFormats()
$ws.Range("A10").ClearFormats()
$ws.Range("A12").ClearFormats()
$ws.Range("A14").ClearFormats()
$ws.Range("A16").ClearFormats()

$ws.Range("A10").Select()
